$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-08-19 Tuesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-08-20 Wednesday", 2)

# Update the division problems in the table. Each problem lives in its own
# table cell, so addressing cells by (row, column) avoids any ambiguity from
# duplicate values that can occur with a pure text search/replace.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "57÷4=" },
    @{ Row = 1;  Col = 2; New = "79÷4=" },
    @{ Row = 1;  Col = 3; New = "12÷6=" },
    @{ Row = 1;  Col = 4; New = "59÷4=" },
    @{ Row = 1;  Col = 5; New = "41÷9=" },

    @{ Row = 5;  Col = 1; New = "53÷6=" },
    @{ Row = 5;  Col = 2; New = "20÷7=" },
    @{ Row = 5;  Col = 3; New = "68÷3=" },
    @{ Row = 5;  Col = 4; New = "37÷2=" },
    @{ Row = 5;  Col = 5; New = "50÷7=" },

    @{ Row = 9;  Col = 1; New = "31÷6=" },
    @{ Row = 9;  Col = 2; New = "96÷7=" },
    @{ Row = 9;  Col = 3; New = "15÷4=" },
    @{ Row = 9;  Col = 4; New = "48÷3=" },
    @{ Row = 9;  Col = 5; New = "74÷9=" },

    @{ Row = 13; Col = 1; New = "76÷9=" },
    @{ Row = 13; Col = 2; New = "81÷5=" },
    @{ Row = 13; Col = 3; New = "50÷5=" },
    @{ Row = 13; Col = 4; New = "71÷3=" },
    @{ Row = 13; Col = 5; New = "89÷2=" },

    @{ Row = 17; Col = 1; New = "50÷5=" },
    @{ Row = 17; Col = 2; New = "23÷4=" },
    @{ Row = 17; Col = 3; New = "76÷3=" },
    @{ Row = 17; Col = 4; New = "25÷9=" },
    @{ Row = 17; Col = 5; New = "84÷4=" }
)

foreach ($u in $updates) {
    $cell = $t.Rows.Item($u.Row).Cells.Item($u.Col)
    $cell.Range.Text = $u.New
}
